$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-10-29 Tuesday" "2024-10-30 Wednesday"

Replace-Text "62×44=" "91×65="
Replace-Text "63×70=" "89×86="
Replace-Text "37×61=" "79×35="
Replace-Text "86×85=" "54×59="
Replace-Text "71×83=" "39×54="

Replace-Text "25×80=" "89×21="
Replace-Text "81×80=" "18×42="
Replace-Text "31×99=" "65×90="
Replace-Text "79×43=" "52×73="
Replace-Text "54×66=" "52×75="

Replace-Text "29×46=" "43×11="
Replace-Text "91×76=" "86×40="
Replace-Text "86×67=" "32×30="
Replace-Text "81×57=" "88×62="
Replace-Text "81×20=" "25×46="

Replace-Text "92×45=" "92×43="
Replace-Text "41×88=" "23×96="
Replace-Text "34×19=" "84×78="
Replace-Text "15×23=" "45×90="
Replace-Text "64×96=" "55×66="

Replace-Text "98×44=" "97×66="
Replace-Text "56×72=" "17×80="
Replace-Text "67×51=" "48×43="
Replace-Text "70×66=" "62×36="
Replace-Text "28×29=" "81×66="
